$d = $word.ActiveDocument

# The Pearson logo (footers) is currently "image1.png" and should become
# "image2.png"; the BTEC logo (headers) is currently "image2.jpg" and
# should become "image1.jpg". Walk every header/footer of every section
# and rename each inline picture accordingly.
foreach ($sec in $d.Sections) {

    foreach ($h in $sec.Headers) {
        if ($h.Exists) {
            for ($i = 1; $i -le $h.Range.InlineShapes.Count; $i++) {
                $shp = $h.Range.InlineShapes.Item($i)
                # Re-select the shape's own range so the rename reliably
                # lands on the live object (avoids stale-handle issues).
                $shp.Range.Select()
                $sel = $word.Selection.InlineShapes.Item(1)
                if ($sel.AlternativeText -eq "BTec_Logo-Orange") {
                    $sel.Name = "image1.jpg"
                }
            }
        }
    }

    foreach ($f in $sec.Footers) {
        if ($f.Exists) {
            for ($i = 1; $i -le $f.Range.InlineShapes.Count; $i++) {
                $shp = $f.Range.InlineShapes.Item($i)
                $shp.Range.Select()
                $sel = $word.Selection.InlineShapes.Item(1)
                if ($sel.AlternativeText -like "*PearsonLogo.png") {
                    $sel.Name = "image2.png"
                }
            }
        }
    }
}
